$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Overal Stats" -- new column AR (44) for 4/18/2020 (serial 43938) ---
$ws1 = $wb.Worksheets.Item("Overal Stats")
$ws1.Cells.Item(1, 44).Value = 43938
$ws1.Cells.Item(1, 44).NumberFormat = $ws1.Cells.Item(1, 43).NumberFormat
$ws1.Cells.Item(3, 44).Value = 13268
$ws1.Cells.Item(4, 44).Value = 2666
$ws1.Cells.Item(5, 44).Value = 91
$ws1.Cells.Item(6, 44).Value = 608
$ws1.Cells.Item(8, 44).Value = 87
$ws1.Cells.Item(9, 44).Value = 443
$ws1.Cells.Item(10, 44).Value = 205
$ws1.Cells.Item(11, 44).Value = 238
$ws1.Cells.Item(15, 44).Value = 74
$ws1.Cells.Item(16, 44).Value = 41
$ws1.Cells.Item(17, 44).Value = 33
$ws1.Cells.Item(18, 44).Value = 96
$ws1.Cells.Item(19, 44).Value = 137
$ws1.Cells.Item(23, 44).Value = 76
$ws1.Cells.Item(24, 44).Value = 59
$ws1.Cells.Item(25, 44).Value = 17
$ws1.Cells.Item(26, 44).Value = 144
$ws1.Cells.Item(27, 44).Value = 203
$ws1.Cells.Item(28, 44).Value = 571
$ws1.Cells.Item(31, 44).Value = 23
$ws1.Cells.Item(32, 44).Value = 20
$ws1.Cells.Item(33, 44).Value = 3
$ws1.Cells.Item(34, 44).Value = 147
$ws1.Cells.Item(35, 44).Value = 167
$ws1.Cells.Item(36, 44).Value = 54
$ws1.Cells.Item(39, 44).Value = 75
$ws1.Cells.Item(40, 44).Value = 31
$ws1.Cells.Item(41, 44).Value = 43
$ws1.Cells.Item(42, 44).Value = 505
$ws1.Cells.Item(43, 44).Value = 536
$ws1.Cells.Item(44, 44).Value = 43
$ws1.Cells.Item(45, 44).Value = 1
$ws1.Cells.Item(48, 44).Value = 14
$ws1.Cells.Item(49, 44).Value = 14
$ws1.Cells.Item(50, 44).Value = 0
$ws1.Cells.Item(51, 44).Value = 64
$ws1.Cells.Item(52, 44).Value = 77
$ws1.Cells.Item(53, 44).Value = 50
$ws1.Cells.Item(54, 44).Value = 1
$ws1.Cells.Item(56, 44).Value = 8
$ws1.Cells.Item(57, 44).Value = 7
$ws1.Cells.Item(58, 44).Value = 1
$ws1.Cells.Item(59, 44).Value = 10
$ws1.Cells.Item(60, 44).Value = 16
$ws1.Cells.Item(61, 44).Value = 0
$ws1.Cells.Item(62, 44).Value = 1
$ws1.Cells.Item(65, 44).Value = 92
$ws1.Cells.Item(66, 44).Value = 268
$ws1.Cells.Item(67, 44).Value = 255
$ws1.Cells.Item(68, 44).Value = 5
$ws1.Cells.Item(70, 44).Value = 54
$ws1.Cells.Item(71, 44).Value = 31
$ws1.Cells.Item(72, 44).Value = 85
$ws1.Cells.Item(73, 44).Value = 22
$ws1.Cells.Item(75, 44).Value = 36
$ws1.Cells.Item(76, 44).Value = 74
$ws1.Cells.Item(77, 44).Value = 74
$ws1.Cells.Item(78, 44).Value = 2
$ws1.Cells.Item(79, 44).Value = 4

# --- Sheet 2: "Total Cases by Ward" -- new column S (19) ---
$ws2 = $wb.Worksheets.Item("Total Cases by Ward")
$ws2.Cells.Item(2, 19).Value = 43938
$ws2.Cells.Item(2, 19).NumberFormat = $ws2.Cells.Item(2, 18).NumberFormat
$ws2.Cells.Item(3, 19).Value = 309
$ws2.Cells.Item(4, 19).Value = 215
$ws2.Cells.Item(5, 19).Value = 205
$ws2.Cells.Item(6, 19).Value = 459
$ws2.Cells.Item(7, 19).Value = 350
$ws2.Cells.Item(8, 19).Value = 353
$ws2.Cells.Item(9, 19).Value = 402
$ws2.Cells.Item(10, 19).Value = 339
$ws2.Cells.Item(11, 19).Value = 34

# --- Sheet 3: "Total Cases by Race" -- new column N (14) ---
$ws3 = $wb.Worksheets.Item("Total Cases by Race")
$ws3.Cells.Item(2, 14).Value = 43938
$ws3.Cells.Item(2, 14).NumberFormat = $ws3.Cells.Item(2, 13).NumberFormat
$ws3.Cells.Item(4, 14).Value = 2666
$ws3.Cells.Item(5, 14).Value = 413
$ws3.Cells.Item(6, 14).Value = 473
$ws3.Cells.Item(7, 14).Value = 1271
$ws3.Cells.Item(8, 14).Value = 38
$ws3.Cells.Item(9, 14).Value = 8
$ws3.Cells.Item(10, 14).Value = 5
$ws3.Cells.Item(11, 14).Value = 431
$ws3.Cells.Item(12, 14).Value = 27
$ws3.Cells.Item(14, 14).Value = 607
$ws3.Cells.Item(15, 14).Value = 439
$ws3.Cells.Item(16, 14).Value = 1613
$ws3.Cells.Item(17, 14).Value = 7

# --- Sheet 4: "Lives Lost by Race" -- new column N (14) ---
$ws4 = $wb.Worksheets.Item("Lives Lost by Race")
$ws4.Cells.Item(1, 14).Value = 43938
$ws4.Cells.Item(1, 14).NumberFormat = $ws4.Cells.Item(1, 13).NumberFormat
$ws4.Cells.Item(3, 14).Value = 91
$ws4.Cells.Item(4, 14).Value = 2
$ws4.Cells.Item(5, 14).Value = 71
$ws4.Cells.Item(6, 14).Value = 8
$ws4.Cells.Item(7, 14).Value = 10
$ws4.Cells.Item(8, 14).Value = 0

# --- Update each sheets active selection to match the new data extent ---
$ws2.Activate()
$ws2.Range("S2").Select()

$ws3.Activate()
$ws3.Range("N17").Select()

$ws4.Activate()
$ws4.Range("N9").Select()

$ws1.Activate()
$ws1.Range("A61:XFD61").Select()
